# Applies the "Added handling of common packages" change:
# the classFields sheet rows describing Customer / CustomerBuilder / OrderServiceImpl
# fields get reordered (their Field Name / Field Type values are shuffled among the
# existing rows), while Class Name / Field Modifier stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# Row -> (Field Name, Field Type)
$updates = @{
    8  = @("amountReserved", "int")
    9  = @("id", "java.lang.Long")
    10 = @("amountAvailable", "int")
    11 = @("name", "java.lang.String")
    12 = @("id", "java.lang.Long")
    13 = @("amountReserved", "int")
    16 = @("SOURCE", "domain.OrderSource")
    17 = @("repository", "com.zatribune.spring.ecommerce.payments.db.repository.CustomerRepository")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
}
